# Apply the "add more elements in the dayahead..." edit.
#
# Summary of changes (per the target diff):
#  - 风力1 (sheet2.xml): scenario-capacity row (row 3, C:F) values shrink,
#    a new hourly index row (row 4, C:Z = 1..24) is added, and a brand new
#    scenario row (row 16, C:Z) is populated. This sheet becomes the
#    workbook's active sheet/tab, with the selection left on I20.
#  - 负荷 (sheet11.xml): selection changes from a single cell (D11) to the
#    header row range B3:Y3 (no data changes there).
#  - 风力8 (sheet9.xml) stops being the active tab (handled automatically
#    by activating 风力1 last).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 负荷: just move the selection to B3:Y3 (no value changes on this sheet)
# ---------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("负荷")
$wsLoad.Activate()
$wsLoad.Range("B3:Y3").Select()

# ---------------------------------------------------------------------
# 风力1: update capacity row, add hour-index row, add new scenario row
# ---------------------------------------------------------------------
$wsWind1 = $wb.Worksheets.Item("风力1")
$wsWind1.Activate()

# Row 3 (C3:F3) capacity figures halved
$wsWind1.Cells.Item(3, 3).Value = 800
$wsWind1.Cells.Item(3, 4).Value = 1600
$wsWind1.Cells.Item(3, 5).Value = 200
$wsWind1.Cells.Item(3, 6).Value = -200

# Row 4 (C4:Z4): new hour-of-day index, 1..24
for ($i = 0; $i -lt 24; $i++) {
    $wsWind1.Cells.Item(4, 3 + $i).Value = $i + 1
}

# Row 16 (C16:Z16): new scenario data row
$row16 = @(1.7, 1.4, 1.2999999999999998, 1.6, 1.6, 1.7000000000000002, 1.9,
            2.1, 2.2999999999999998, 2.1, 2.1, 2, 1.9, 2, 1.7999999999999998,
            1.9, 2.1, 3, 2.5, 2.2999999999999998, 2.2000000000000002, 2.5,
            2.2000000000000002, 1.7000000000000002)
for ($i = 0; $i -lt 24; $i++) {
    $wsWind1.Cells.Item(16, 3 + $i).Value = $row16[$i]
}

# Finish with the selection/active-tab state the diff records: 风力1 is the
# active sheet/tab and I20 is selected.
$wsWind1.Range("I20").Select()
